$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose content is plain text (coin name / link columns) - safe to assign as-is.
$plainUpdates = @{
    'B45' = 'EnergySwap'
    'C45' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'B46' = 'VeChain'
    'C46' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
}

foreach ($ref in $plainUpdates.Keys) {
    $ws.Range($ref).Value = $plainUpdates[$ref]
}

# Cells whose content looks numeric (price/volume columns) but must stay text, matching
# the sheet's existing inline-string cells. A leading apostrophe forces text entry (as it
# would in the Excel UI) and resetting the style afterwards avoids leaving a stray
# Text-format (NumberFormat "@") style behind on the cell.
$textUpdates = [ordered]@{
    'D2'  = '39.874.23';   'E2'  = '  -0.55%  '
    'D3'  = '2.205.07';    'E3'  = '  -1.39%  '
    'E4'  = '  +0.03%  '
    'D5'  = '293.23';      'E5'  = '  -0.58%  '
    'D6'  = '86.72';       'E6'  = '  +0.26%  '
    'E7'  = '  -1.43%  '
    'E8'  = '  +0.07%  '
    'D9'  = '0.472';       'E9'  = '  -0.21%  '
    'E10' = '  -2.56%  '
    'D11' = '29.70';       'E11' = '  -4.69%  '
    'D12' = '49.06';       'E12' = '  +4.31%  '
    'E14' = '  -0.10%  '
    'D15' = '2.548.96'
    'D16' = '13.66';       'E16' = '  -3.65%  '
    'D17' = '2.204.75';    'E17' = '  +1.15%  '
    'D18' = '0.725';       'E18' = '  -0.68%  '
    'D19' = '39.789.55';   'E19' = '  -0.51%  '
    'D20' = '0.0₃0881';    'E20' = '  -1.05%  '
    'D21' = '11.26';       'E21' = '  +3.00%  '
    'E22' = '  -1.04%  '
    'D23' = '65.05';       'E23' = '  -0.68%  '
    'D24' = '235.70';      'E24' = '  +0.16%  '
    'E25' = '  -0.03%  '
    'E26' = '  -0.61%  '
    'E27' = '  -2.74%  '
    'D28' = '22.39';       'E28' = '  -2.03%  '
    'E29' = '  -3.65%  '
    'D30' = '9.13';        'E30' = '  -1.19%  '
    'D31' = '155.02';      'E31' = '  +1.70%  '
    'D32' = '31.61';       'E32' = '  -5.48%  '
    'E33' = '  -0.01%  '
    'D34' = '4.87';        'E34' = '  -0.44%  '
    'D35' = '0.0709';      'E35' = '  -1.58%  '
    'E36' = '  -1.61%  '
    'E37' = '  +3.51%  '
    'E38' = '  -0.44%  '
    'D39' = '0.0975';      'E39' = '  -3.02%  '
    'D40' = '15.42';       'E40' = '  -6.12%  '
    'E41' = '  -2.52%  '
    'D42' = '2.119.74';    'E42' = '  +3.59%  '
    'D43' = '3.71';        'E43' = '  -3.42%  '
    'D44' = '2.10';        'E44' = '  -5.96%  '
    'D45' = '17.78';       'E45' = '  +8.86%  '
    'D46' = '0.0266';      'E46' = '  -1.56%  '
    'D47' = '9.61';        'E47' = '  -3.93%  '
    'E48' = '  +3.44%  '
    'D49' = '2.416.43';    'E49' = '  -1.42%  '
    'E50' = '  -0.86%  '
    'E51' = '  -0.18%  '
}

foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = "'" + $textUpdates[$ref]
    $ws.Range($ref).Style = 'Normal'
}
